# Recognize individual words in the "Splate" (abbreviation) column:
#  - "DFA,DA"  -> "DFA, DA"   (add a space after the comma)
#  - "Berlin"  -> "Berlin, TXL" (append the related abbreviation)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the two text values (shared strings get updated accordingly).
$ws.Range("A6").Value = "DFA, DA"
$ws.Range("A7").Value = "Berlin, TXL"

# The longer text in column A required the column to be widened.
# (211/6 is the exact midpoint of the ColumnWidth bucket that rounds to the
# stored width closest to the author's target width.)
$ws.Columns.Item(1).ColumnWidth = 35.166666666666664

# The author's active/selected cell ended up on the newly edited A7 cell.
$ws.Range("A7").Select() | Out-Null
